$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that get new data in columns J, K, L (value 0).
# Rows 6, 12, 18 are "special" shorter rows (Switch/Station) that are
# skipped - only their row span metadata changes, not their cell data.
$dataRows = @(1,2,3,4,5,7,8,9,10,11,13,14,15,16,17)

foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 10).Value = 0   # J
    $ws.Cells.Item($r, 11).Value = 0   # K
    $ws.Cells.Item($r, 12).Value = 0   # L
}
